$d = $word.ActiveDocument

# Delete the 5 "group list" paragraphs (Bhavjot Pal..., Yiyuan Dong..., Huu Minh...,
# Xin Zhao..., Dennis Audu...) leaving the title and the two blank paragraphs
# around them intact, producing a blank groups list.
$start = $d.Paragraphs.Item(3).Range.Start
$end = $d.Paragraphs.Item(7).Range.End

$r = $d.Range($start, $end)
$r.Delete()
